$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("H1").Value = "reference_period"
$ws.Range("I1").Value = "remarks"

# New data cells
$ws.Range("H2").Value = 2020
$ws.Range("I2").Value = "Test note"

# Column H width (diff: width="14.59")
$ws.Columns("H").ColumnWidth = 13.75

# Re-apply the Normal style on the original data range -> marks the cellXf
# as an explicit (applyFont) style distinct from the default, matching the
# font/format refresh the original data range received.
$ws.Range("A1:G2").Style = "Normal"

# Selection moves to F8 per diff
$ws.Range("F8").Select()
